$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 174: 2025-11-26, 四方坪站充电量(kw)
$ws.Range("A174").Value = 45987
$ws.Range("B174").Value = "四方坪站充电量(kw)"
$ws.Range("C174").Value = 511.36600000000004
$ws.Range("D174").Value = 847.13100000000009
$ws.Range("E174").Value = 339.85699999999997
$ws.Range("F174").Value = 352.62099999999998
$ws.Range("G174").Value = 341.41800000000001
$ws.Range("H174").Value = 628.4910000000001
$ws.Range("I174").Value = 321.53799999999995
$ws.Range("J174").Value = 172.02900000000002
$ws.Range("K174").Value = 182.547
$ws.Range("L174").Value = 156.30199999999999
$ws.Range("M174").Value = 226.16000000000003
$ws.Range("N174").Value = 193.85900000000001
$ws.Range("O174").Value = 575.49
$ws.Range("P174").Value = 1124.0030000000002
$ws.Range("Q174").Value = 630.85900000000004
$ws.Range("R174").Value = 493.69799999999992
$ws.Range("S174").Value = 199.21999999999997
$ws.Range("T174").Value = 125.39899999999999
$ws.Range("U174").Value = 182.98899999999998
$ws.Range("V174").Value = 84.725999999999999
$ws.Range("W174").Value = 171.34
$ws.Range("X174").Value = 79.619
$ws.Range("Y174").Value = 17.3
$ws.Range("Z174").Value = 20.420000000000002

# Row 175: 2025-11-26, 高岭站充电量(kw)
$ws.Range("A175").Value = 45987
$ws.Range("B175").Value = "高岭站充电量(kw)"
$ws.Range("C175").Value = 457.23899999999998
$ws.Range("D175").Value = 738.03099999999984
$ws.Range("E175").Value = 307.92500000000001
$ws.Range("F175").Value = 58.226999999999997
$ws.Range("G175").Value = 0
$ws.Range("H175").Value = 81.893000000000001
$ws.Range("I175").Value = 311.21799999999996
$ws.Range("J175").Value = 186.70699999999999
$ws.Range("K175").Value = 173.45599999999999
$ws.Range("L175").Value = 320.85599999999999
$ws.Range("M175").Value = 292.75100000000003
$ws.Range("N175").Value = 364.53900000000004
$ws.Range("O175").Value = 592.99099999999999
$ws.Range("P175").Value = 487.63200000000006
$ws.Range("Q175").Value = 113.9
$ws.Range("R175").Value = 168.45699999999999
$ws.Range("S175").Value = 244.53
$ws.Range("T175").Value = 53.773000000000003
$ws.Range("U175").Value = 80.527999999999992
$ws.Range("V175").Value = 70.796999999999997
$ws.Range("W175").Value = 26.066000000000003
$ws.Range("X175").Value = 26.277000000000001
$ws.Range("Y175").Value = 29.516000000000002
$ws.Range("Z175").Value = 76.915999999999997

# Update selection / scroll to mirror the author's final view state
$ws.Range("E179").Select() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollRow = 169
$win.ScrollColumn = 1
